$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.98"
$ws.Range("E2").Value = "'-0.67%"
$ws.Range("D3").Value = "'26.66"
$ws.Range("E3").Value = "'4.63%"
$ws.Range("D4").Value = "'5.159"
$ws.Range("E4").Value = "'0.32%"
$ws.Range("D5").Value = "'0.05619"
$ws.Range("E5").Value = "'0.49%"
$ws.Range("D6").Value = "'6.492"
$ws.Range("E6").Value = "'0.05%"
$ws.Range("D7").Value = "'0.8177"
$ws.Range("E7").Value = "'-0.01%"
$ws.Range("D8").Value = "'0.8316"
$ws.Range("E8").Value = "'-1.93%"
$ws.Range("D9").Value = "'0.1325"
$ws.Range("E9").Value = "'-1.13%"
$ws.Range("D10").Value = "'0.06934"
$ws.Range("E10").Value = "'-0.21%"
$ws.Range("D11").Value = "'0.02899"
$ws.Range("E11").Value = "'1.44%"
$ws.Range("D12").Value = "'0.09382"
$ws.Range("E12").Value = "'-0.16%"
$ws.Range("D13").Value = "'0.001518"
$ws.Range("E13").Value = "'-0.38%"
$ws.Range("D14").Value = "'0.0005959"
$ws.Range("E14").Value = "'-93.89%"
$ws.Range("D15").Value = "'0.006218"
$ws.Range("E15").Value = "'-0.11%"
$ws.Range("D16").Value = "'3.644"
$ws.Range("E16").Value = "'3.17%"
$ws.Range("D17").Value = "'3.024"
$ws.Range("E17").Value = "'-0.12%"
$ws.Range("D18").Value = "'2.191"
$ws.Range("E18").Value = "'3.45%"
$ws.Range("E19").Value = "'-1.78%"
$ws.Range("D20").Value = "'0.03095"
$ws.Range("E20").Value = "'-4.42%"
$ws.Range("D21").Value = "'0.1291"
$ws.Range("E21").Value = "'-2.12%"
$ws.Range("D22").Value = "'3.747"
$ws.Range("E22").Value = "'0.00%"
$ws.Range("D23").Value = "'0.04584"
$ws.Range("E23").Value = "'-2.19%"
$ws.Range("E24").Value = "'-2.43%"
$ws.Range("E25").Value = "'-1.58%"
$ws.Range("E26").Value = "'-2.47%"
$ws.Range("E27").Value = "'2.12%"
$ws.Range("E28").Value = "'0.70%"
$ws.Range("D41").Value = "'0.006160"
$ws.Range("E41").Value = "'0.69%"
$ws.Range("D42").Value = "'0.1052"
$ws.Range("E42").Value = "'0.00%"
$ws.Range("D43").Value = "'0.002589"
$ws.Range("E43").Value = "'4.80%"
$ws.Range("D44").Value = "'0.008095"
$ws.Range("E44").Value = "'4.13%"
$ws.Range("D45").Value = "'0.00005348"
$ws.Range("E45").Value = "'0.55%"
$ws.Range("E46").Value = "'0.04%"
$ws.Range("E47").Value = "'-18.33%"
$ws.Range("D48").Value = "'0.002587"
$ws.Range("E48").Value = "'21.68%"
$ws.Range("E49").Value = "'0.04%"
$ws.Range("E50").Value = "'0.04%"
